# RF011 - Gerenciar Competencias (Portfolio): grammatical fix from
# masculine ("o Competencias (Portfolio)" / "dos Competencias (Portfolio)")
# to feminine ("a Competencia (Portfolio)" / "das Competencias (Portfolio)")
# agreement, matching "Competencia" being a feminine noun in Portuguese.
# (The commit's "1.0 to 1.1" note refers to the authoring tool/template
# version, not the workbook's own "Version:" cell, which the OOXML diff
# leaves untouched.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "SYSTEM exibe a listagem dos Competencias (Portfolio) cadastrados apenas
#  para visualizacao com a opcao 'Ajuda'"
#  -> "...listagem das Competencias (Portfolio) cadastradas apenas..."
$refs1 = @("D10")
$text1 = "SYSTEM exibe a listagem das Competencias (Portfolio) cadastradas apenas para visualizacao com a opcao 'Ajuda'"
foreach ($ref in $refs1) {
    $ws.Range($ref).Value = $text1
}

# "SYSTEM exibe a listagem dos Competencias (Portfolio) cadastrados com
#  opcoes de 'Novo', 'Editar', 'Excluir' e 'Ajuda'"
#  -> "...listagem das Competencias (Portfolio) cadastradas com..."
$refs2 = @("D17","D27","D37","D47","D62","D77","D93","D107","D121","D136","D151","D166","D182","D196","D210","D225","D239","D253","D268","D281","D294","D308","D322","D336","D351","D364","D377")
$text2 = "SYSTEM exibe a listagem das Competencias (Portfolio) cadastradas com opcoes de 'Novo', 'Editar', 'Excluir' e 'Ajuda'"
foreach ($ref in $refs2) {
    $ws.Range($ref).Value = $text2
}

# "SYSTEM destaca o Competencias (Portfolio) selecionado na listagem"
#  -> "SYSTEM destaca a Competencia (Portfolio) selecionada na listagem"
$refs3 = @("D18","D28","D38","D48","D63","D78","D94","D108","D122","D137","D152","D167","D183","D197","D211")
$text3 = "SYSTEM destaca a Competencia (Portfolio) selecionada na listagem"
foreach ($ref in $refs3) {
    $ws.Range($ref).Value = $text3
}

# "Lider de Pessoas clica na opcao 'Excluir' para excluir o Competencias
#  (Portfolio) selecionado"
#  -> "...para excluir a Competencia (Portfolio) selecionada"
$refs4 = @("B19","B29","B39")
$text4 = "Lider de Pessoas clica na opcao 'Excluir' para excluir a Competencia (Portfolio) selecionada"
foreach ($ref in $refs4) {
    $ws.Range($ref).Value = $text4
}

# "SYSTEM exibe a listagem dos Competencias (Portfolio) com o Competencias
#  (Portfolio) excluido"
#  -> "...listagem das Competencias (Portfolio) com a Competencia (Portfolio) excluida"
$refs5 = @("D20")
$text5 = "SYSTEM exibe a listagem das Competencias (Portfolio) com a Competencia (Portfolio) excluida"
foreach ($ref in $refs5) {
    $ws.Range($ref).Value = $text5
}

# "SYSTEM exibe a listagem dos Competencias (Portfolio) sem o Competencias
#  (Portfolio) excluido"
#  -> "...listagem das Competencias (Portfolio) sem a Competencia (Portfolio) excluida"
$refs6 = @("D30")
$text6 = "SYSTEM exibe a listagem das Competencias (Portfolio) sem a Competencia (Portfolio) excluida"
foreach ($ref in $refs6) {
    $ws.Range($ref).Value = $text6
}

# "Lider de Pessoas clica na opcao 'Editar' para modificar o Competencias
#  (Portfolio) selecionado"
#  -> "...para modificar a Competencia (Portfolio) selecionada"
$refs7 = @("B49","B64","B79","B95","B109","B123","B138","B153","B168","B184","B198","B212")
$text7 = "Lider de Pessoas clica na opcao 'Editar' para modificar a Competencia (Portfolio) selecionada"
foreach ($ref in $refs7) {
    $ws.Range($ref).Value = $text7
}
